# Updated the gear ratios for the comp bot (docs/motormapping.xlsx, Sheet1).
#
# Sheet1 lists, for each mechanism, a "Gear Ratio" input cell (column D)
# whose ticks-per-degree figures (column F/G) are derived by formula, so
# only the raw gear-ratio inputs need to change - the dependent formulas
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shooter gear ratio: 101.56 -> 250
$ws.Range("D23").Value = 250

# Elevator gear ratio: 40 -> 29.33
$ws.Range("D26").Value = 29.33

# Pivot gear ratio: 12.5 -> 14.81
$ws.Range("D29").Value = 14.81

# Leave the cursor where the author last left it on Sheet1 (was G21).
$ws.Range("H30").Select()
